# "Added some openvpn lines to search for"
# Adds two new rows (OpenVPN comp-lzo related log messages) to the
# "IPSec" worksheet, right after the existing log-message table, and
# makes that sheet the active/selected one (mirroring the author's
# Excel session state when the file was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPSec")

# Bring the IPSec sheet to the front, matching the saved workbook view
# (tabSelected moves from "Connectivity+Modem" to "IPSec").
$ws.Activate()

# Fill the two new rows. The cells are written in the same order the
# original author entered them (Message column first, then the rest),
# which is what determines the insertion order of the new shared-string
# table entries.
$ws.Range("C9").Value = "'comp-lzo' is present in remote config but missing in local config, remote='comp-lzo'"
$ws.Range("B9").Value = "OPENVPN"
$ws.Range("D9").Value = "Comp-lzo is a compression option for OpenVPN.  Cradlepoint routers do not have a way to turn on Comp-lzo.  To get the OpenVPN tunnel to connect, the remote side of the tunnel will have to turn off the comp-lzo flag in their OpenVPN config"

$ws.Range("C10").Value = "write to TUN/TAP : Invalid argument \(code=22\)"
$ws.Range("A10").Value = "ERR"
$ws.Range("D10").Value = "The remote side of the OpenVPN tunnel is trying to push the comp-lzo option to the Cradlepoint.  Cradlepoints do not support using Comp-lzo, so the tunnel wont function properly.  To get the OpenVPN tunnel to connect, the remote side of the tunnel will have to turn off the comp-lzo flag in their OpenVPN config"

$ws.Range("A9").Value = "WARN"
$ws.Range("B10").Value = "OPENVPN"

# Row heights (wrapped text needs more vertical room than the default).
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 60

# Leave the selection where the author left it.
$ws.Range("A11").Select()
